$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180 - this shifts the existing rows
# 180-189 down to 181-190, matching the diff (which is effectively a
# single new weekly record inserted in date order, pushing the rest of
# the series down by one row).
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new weekly record.
$ws.Range("A180").Value = 7
$ws.Range("B180").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C180").Value = "Ñuble"
$ws.Range("D180").Value = 44568
$ws.Range("E180").Value = 16
$ws.Range("F180").Value = "Fruta"
$ws.Range("G180").Value = 100108
$ws.Range("H180").Value = "Tropicales y subtropicales"
$ws.Range("I180").Value = 100108005
$ws.Range("J180").Value = "Piña"
$ws.Range("K180").Value = "Caramelo"
$ws.Range("L180").Value = "Segunda"
$ws.Range("M180").Value = 120
$ws.Range("N180").Value = 16000
$ws.Range("O180").Value = 17000
$ws.Range("P180").Value = 16500
$ws.Range("Q180").Value = "`$/caja 14 unidades"
$ws.Range("R180").Value = "Ecuador"
$ws.Range("S180").Value = 1179
$ws.Range("T180").Value = 14
